$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1829,1).Value = 4
$ws.Cells.Item(1829,2).Value = 3
$ws.Cells.Item(1829,3).Value = 5
$ws.Cells.Item(1829,4).Value = 0

$ws.Cells.Item(1830,1).Value = 6
$ws.Cells.Item(1830,2).Value = 2
$ws.Cells.Item(1830,3).Value = 5
$ws.Cells.Item(1830,4).Value = 0

$ws.Cells.Item(1831,1).Value = 3
$ws.Cells.Item(1831,2).Value = 3
$ws.Cells.Item(1831,3).Value = 4
$ws.Cells.Item(1831,4).Value = 0

$ws.Cells.Item(1832,1).Value = 5
$ws.Cells.Item(1832,2).Value = 2
$ws.Cells.Item(1832,3).Value = 5
$ws.Cells.Item(1832,4).Value = 0

$ws.Cells.Item(1833,1).Value = 5
$ws.Cells.Item(1833,2).Value = 0
$ws.Cells.Item(1833,3).Value = 2
$ws.Cells.Item(1833,4).Value = 3

$ws.Cells.Item(1834,1).Value = 5
$ws.Cells.Item(1834,2).Value = 3
$ws.Cells.Item(1834,3).Value = 4
$ws.Cells.Item(1834,4).Value = 0

$ws.Cells.Item(1835,1).Value = 4
$ws.Cells.Item(1835,2).Value = 2
$ws.Cells.Item(1835,3).Value = 3
$ws.Cells.Item(1835,4).Value = 1

$ws.Cells.Item(1836,1).Value = 3
$ws.Cells.Item(1836,2).Value = 0
$ws.Cells.Item(1836,3).Value = 2
$ws.Cells.Item(1836,4).Value = 2

$ws.Cells.Item(1837,1).Value = 6
$ws.Cells.Item(1837,2).Value = 2
$ws.Cells.Item(1837,3).Value = 2
$ws.Cells.Item(1837,4).Value = 1

$ws.Cells.Item(1838,1).Value = 6
$ws.Cells.Item(1838,2).Value = 0
$ws.Cells.Item(1838,3).Value = 6
$ws.Cells.Item(1838,4).Value = 3

$ws.Cells.Item(1839,1).Value = 3
$ws.Cells.Item(1839,2).Value = 0
$ws.Cells.Item(1839,3).Value = 3
$ws.Cells.Item(1839,4).Value = 3

$ws.Cells.Item(1840,1).Value = 5
$ws.Cells.Item(1840,2).Value = 2
$ws.Cells.Item(1840,3).Value = 5
$ws.Cells.Item(1840,4).Value = 0

$ws.Cells.Item(1841,1).Value = 4
$ws.Cells.Item(1841,2).Value = 0
$ws.Cells.Item(1841,3).Value = 3
$ws.Cells.Item(1841,4).Value = 2

$ws.Cells.Item(1842,1).Value = 5
$ws.Cells.Item(1842,2).Value = 1
$ws.Cells.Item(1842,3).Value = 5
$ws.Cells.Item(1842,4).Value = 2

$ws.Cells.Item(1843,1).Value = 4
$ws.Cells.Item(1843,2).Value = 1
$ws.Cells.Item(1843,3).Value = 4
$ws.Cells.Item(1843,4).Value = 2

$ws.Cells.Item(1844,1).Value = 4
$ws.Cells.Item(1844,2).Value = 0
$ws.Cells.Item(1844,3).Value = 3
$ws.Cells.Item(1844,4).Value = 3

$ws.Cells.Item(1845,1).Value = 5
$ws.Cells.Item(1845,2).Value = 2
$ws.Cells.Item(1845,3).Value = 5
$ws.Cells.Item(1845,4).Value = 0

$ws.Cells.Item(1846,1).Value = 5
$ws.Cells.Item(1846,2).Value = 0
$ws.Cells.Item(1846,3).Value = 6
$ws.Cells.Item(1846,4).Value = 3

$ws.Cells.Item(1847,1).Value = 5
$ws.Cells.Item(1847,2).Value = 0
$ws.Cells.Item(1847,3).Value = 5
$ws.Cells.Item(1847,4).Value = 2

$ws.Cells.Item(1848,1).Value = 5
$ws.Cells.Item(1848,2).Value = 2
$ws.Cells.Item(1848,3).Value = 6
$ws.Cells.Item(1848,4).Value = 0

$ws.Cells.Item(1849,1).Value = 5
$ws.Cells.Item(1849,2).Value = 3
$ws.Cells.Item(1849,3).Value = 2
$ws.Cells.Item(1849,4).Value = 0

$ws.Cells.Item(1850,1).Value = 3
$ws.Cells.Item(1850,2).Value = 0
$ws.Cells.Item(1850,3).Value = 5
$ws.Cells.Item(1850,4).Value = 3

$ws.Cells.Item(1851,1).Value = 4
$ws.Cells.Item(1851,2).Value = 2
$ws.Cells.Item(1851,3).Value = 6
$ws.Cells.Item(1851,4).Value = 1

$ws.Cells.Item(1852,1).Value = 3
$ws.Cells.Item(1852,2).Value = 3
$ws.Cells.Item(1852,3).Value = 2
$ws.Cells.Item(1852,4).Value = 0

$ws.Range("A1853").Select()
